$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-18 Sunday", "2026-01-19 Monday"),
    @("59÷9=", "76÷5="),
    @("39÷4=", "26÷7="),
    @("89÷7=", "26÷5="),
    @("19÷5=", "65÷5="),
    @("21÷6=", "20÷7="),
    @("63÷4=", "12÷4="),
    @("78÷5=", "79÷4="),
    @("61÷9=", "97÷8="),
    @("57÷5=", "21÷8="),
    @("16÷2=", "90÷3="),
    @("76÷2=", "74÷8="),
    @("34÷2=", "72÷2="),
    @("33÷2=", "46÷3="),
    @("93÷6=", "63÷3="),
    @("83÷8=", "33÷6="),
    @("28÷5=", "69÷4="),
    @("96÷8=", "75÷3="),
    @("75÷4=", "98÷7="),
    @("34÷7=", "20÷3="),
    @("90÷2=", "24÷7="),
    @("83÷3=", "32÷8="),
    @("87÷9=", "41÷6="),
    @("22÷2=", "95÷7="),
    @("35÷5=", "96÷7="),
    @("14÷3=", "82÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
